$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Developer name
$ws.Range("C3").Value = "Jashanpreet Singh Sidhu"

# Preconditions (E), Method Inputs (F), Expected Result (G) for test rows 7-22
$ws.Range("E7").Value = "None"
$ws.Range("F7").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("G7").Value = "Attributes are set "

$ws.Range("E8").Value = "None"
$ws.Range("F8").Value = "account_number = 200`nclient_number = 200`nbalance = ""two hundred"""
$ws.Range("G8").Value = "balance attribute set to 0"

$ws.Range("E9").Value = "None"
$ws.Range("F9").Value = "account_number = ""two hundred""`nclient_number = 200`nbalance = 200"
$ws.Range("G9").Value = "ValueError "

$ws.Range("E10").Value = "None"
$ws.Range("F10").Value = "account_number = 200`nclient_number = ""two hundred""`nbalance = 200"
$ws.Range("G10").Value = "ValueError "

$ws.Range("E11").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "account_number"

$ws.Range("E12").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F12").Value = "None"
$ws.Range("G12").Value = "client_number"

$ws.Range("E13").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F13").Value = "None"
$ws.Range("G13").Value = "balance attribute set to 0"

$ws.Range("E14").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F14").Value = "amount = 100"
$ws.Range("G14").Value = "balance attributes updates"

$ws.Range("E15").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F15").Value = "amount = -100"
$ws.Range("G15").Value = "balance attributes updates"

$ws.Range("E16").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F16").Value = "amount = ""hundred"""
$ws.Range("G16").Value = "balance attribute unchanged"

$ws.Range("E17").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F17").Value = "amount = 100"
$ws.Range("G17").Value = "balance attributes updates"

$ws.Range("E18").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F18").Value = "amount = -100"
$ws.Range("G18").Value = "ValueError"

$ws.Range("E19").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F19").Value = "amount = 100"
$ws.Range("G19").Value = "balance attributes updates"

$ws.Range("E20").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F20").Value = "amount = -100"
$ws.Range("G20").Value = "ValueError"

$ws.Range("E21").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F21").Value = "amount = 400"
$ws.Range("G21").Value = "ValueError when amount"

$ws.Range("E22").Value = "account_number = 200`nclient_number = 200`nbalance = 200"
$ws.Range("F22").Value = "None"
$ws.Range("G22").Value = "Account Number: 200 Balance: `$200.00"

# View settings
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Range("F8").Select()
